$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose numeric-looking values would otherwise be
# auto-converted to numbers by Excel, to preserve the original text cell type.
$textCells = @("D5","D6","D7","D12","D15","D18","D19","D21","D22","D25","D26","D27","D28","D29","D32","D34","D38","D40","D41","D42","D45","D46","D49","D51")
foreach ($c in $textCells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range('D2').Value = '37.145.82'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').Value = '2.051.20'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '249.76'
$ws.Range('E5').Value = '  +0.30%  '
$ws.Range('D6').Value = '0.665'
$ws.Range('D7').Value = '59.33'
$ws.Range('E7').Value = '  +6.87%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +0.75%  '
$ws.Range('E10').Value = '  -1.82%  '
$ws.Range('E11').Value = '  +1.70%  '
$ws.Range('D12').Value = '16.14'
$ws.Range('E12').Value = '  +6.43%  '
$ws.Range('D13').Value = '2.354.61'
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('E14').Value = '  -0.55%  '
$ws.Range('D15').Value = '5.59'
$ws.Range('E15').Value = '  +6.33%  '
$ws.Range('D16').Value = '2.053.76'
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('D17').Value = '37.120.84'
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').Value = '17.47'
$ws.Range('E18').Value = '  +22.86%  '
$ws.Range('D19').Value = '74.84'
$ws.Range('E19').Value = '  +3.34%  '
$ws.Range('D20').Value = '0.0₃0902'
$ws.Range('E20').Value = '  -1.05%  '
$ws.Range('D21').Value = '5.37'
$ws.Range('E21').Value = '  +0.52%  '
$ws.Range('D22').Value = '237.43'
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('E24').Value = '  -1.10%  '
$ws.Range('D25').Value = '2.20'
$ws.Range('E25').Value = '  +10.89%  '
$ws.Range('D26').Value = '168.52'
$ws.Range('E26').Value = '  -0.89%  '
$ws.Range('D27').Value = '9.35'
$ws.Range('E27').Value = '  +3.27%  '
$ws.Range('D28').Value = '19.94'
$ws.Range('E28').Value = '  -1.29%  '
$ws.Range('D29').Value = '0.125'
$ws.Range('E29').Value = '  +1.23%  '
$ws.Range('E30').Value = '  +7.98%  '
$ws.Range('E31').Value = '  +4.28%  '
$ws.Range('D32').Value = '0.0616'
$ws.Range('E32').Value = '  -1.47%  '
$ws.Range('E33').Value = '  +4.17%  '
$ws.Range('D34').Value = '0.0899'
$ws.Range('E34').Value = '  +3.43%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('E36').Value = '  -1.26%  '
$ws.Range('E37').Value = '  -2.23%  '
$ws.Range('D38').Value = '0.110'
$ws.Range('E38').Value = '  +6.46%  '
$ws.Range('E39').Value = '  +0.95%  '
$ws.Range('B40').Value = 'THORChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D40').Value = '5.31'
$ws.Range('E40').Value = '  +32.28%  '
$ws.Range('B41').Value = 'HuobiToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D41').Value = '3.16'
$ws.Range('E41').Value = '  +13.41%  '
$ws.Range('D42').Value = '17.61'
$ws.Range('E42').Value = '  -3.02%  '
$ws.Range('E44').Value = '  -1.31%  '
$ws.Range('D45').Value = '96.20'
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('D46').Value = '2.47'
$ws.Range('E46').Value = '  +1.88%  '
$ws.Range('D47').Value = '1.284.91'
$ws.Range('E47').Value = '  -1.06%  '
$ws.Range('E48').Value = '  -1.23%  '
$ws.Range('D49').Value = '6.80'
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('D50').Value = '2.247.26'
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('D51').Value = '3.42'
$ws.Range('E51').Value = '  -21.33%  '

# Restore default style on the forced-text cells so formatting matches the original.
foreach ($c in $textCells) { $ws.Range($c).Style = "Normal" }
